$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 'Multiline Retail(7)'
$ws.Range("B3").Value = 0.637158856603528
$ws.Range("A4").Value = 'Multi-Utilities(18)'
$ws.Range("B4").Value = 0.5909373411224138
$ws.Range("A5").Value = 'Containers & Packaging(12)'
$ws.Range("B5").Value = 0.5820552138888447
$ws.Range("A6").Value = 'Machinery(85)'
$ws.Range("B6").Value = 0.5773156141689123
$ws.Range("A7").Value = 'Road & Rail(22)'
$ws.Range("B7").Value = 0.5710044321418385
$ws.Range("A8").Value = 'Auto Components(21)'
$ws.Range("B8").Value = 0.5537157241117897
$ws.Range("A9").Value = 'Air Freight & Logistics(11)'
$ws.Range("B9").Value = 0.5118791143994464
$ws.Range("A10").Value = 'Energy Equipment & Services(32)'
$ws.Range("B10").Value = 0.5006236236380753
$ws.Range("A11").Value = 'Electric Utilities(28)'
$ws.Range("B11").Value = 0.5001381885300904
$ws.Range("A12").Value = 'Airlines(14)'
$ws.Range("B12").Value = 0.4917499218936111
$ws.Range("A13").Value = 'Specialty Retail(58)'
$ws.Range("B13").Value = 0.4879427848495779
$ws.Range("A14").Value = 'Trading Companies & Distributors(25)'
$ws.Range("B14").Value = 0.4854130637614995
$ws.Range("A15").Value = 'Gas Utilities(12)'
$ws.Range("B15").Value = 0.4786730779818642
$ws.Range("A16").Value = 'Wireless Telecommunication Services(14)'
$ws.Range("B16").Value = 0.4702411239861925
$ws.Range("A17").Value = 'Hotels, Restaurants & Leisure(50)'
$ws.Range("B17").Value = 0.4685900344946026
$ws.Range("A18").Value = 'Chemicals(51)'
$ws.Range("B18").Value = 0.4634824288800502
$ws.Range("A19").Value = 'Media(42)'
$ws.Range("B19").Value = 0.4632595642572921
$ws.Range("A20").Value = 'Construction & Engineering(20)'
$ws.Range("B20").Value = 0.4551704769440679
$ws.Range("A21").Value = 'Semiconductors & Semiconductor Equipment(68)'
$ws.Range("B21").Value = 0.4341880161375474
$ws.Range("A22").Value = 'Leisure Products(11)'
$ws.Range("B22").Value = 0.4273113347199632
$ws.Range("A23").Value = 'Building Products(23)'
$ws.Range("B23").Value = 0.4193499308950515
$ws.Range("A24").Value = 'Household Durables(39)'
$ws.Range("B24").Value = 0.4173637377507899
$ws.Range("A25").Value = 'Capital Markets(75)'
$ws.Range("B25").Value = 0.4068347560461073
$ws.Range("A26").Value = 'Household Products(10)'
$ws.Range("B26").Value = 0.4064015974987705
$ws.Range("A27").Value = 'Marine(15)'
$ws.Range("B27").Value = 0.4033073324332358
$ws.Range("A28").Value = 'Aerospace & Defense(37)'
$ws.Range("B28").Value = 0.3975633353773633
$ws.Range("A29").Value = 'Communications Equipment(45)'
$ws.Range("B29").Value = 0.3803307251461573
$ws.Range("A30").Value = 'Oil, Gas & Consumable Fuels(122)'
$ws.Range("B30").Value = 0.3755636762179344
$ws.Range("A31").Value = 'Insurance(75)'
$ws.Range("B31").Value = 0.3678833024531595
$ws.Range("A32").Value = 'Commercial Services & Supplies(52)'
$ws.Range("B32").Value = 0.3607899127967236
$ws.Range("A33").Value = 'Consumer Finance(15)'
$ws.Range("B33").Value = 0.3521445803575914
$ws.Range("A34").Value = 'Textiles, Apparel & Luxury Goods(29)'
$ws.Range("B34").Value = 0.3515114243143465
$ws.Range("A35").Value = 'Diversified Consumer Services(17)'
$ws.Range("B35").Value = 0.3499089347939527
$ws.Range("A36").Value = 'Metals & Mining(89)'
$ws.Range("B36").Value = 0.3388671756232035
$ws.Range("A37").Value = 'Water Utilities(12)'
$ws.Range("B37").Value = 0.3368986889980476
$ws.Range("A38").Value = 'Diversified Telecommunication Services(20)'
$ws.Range("B38").Value = 0.3291281176173813
$ws.Range("A39").Value = 'Professional Services(35)'
$ws.Range("B39").Value = 0.3241847688660761
$ws.Range("A40").Value = 'Electrical Equipment(28)'
$ws.Range("B40").Value = 0.3177397643706062
$ws.Range("A41").Value = 'Life Sciences Tools & Services(19)'
$ws.Range("B41").Value = 0.3115688160483179
$ws.Range("A42").Value = 'Banks(246)'
$ws.Range("B42").Value = 0.3029498278140831
$ws.Range("A43").Value = 'Food & Staples Retailing(15)'
$ws.Range("B43").Value = 0.3026795069584989
$ws.Range("A44").Value = 'Software(66)'
$ws.Range("B44").Value = 0.2980981953635595
$ws.Range("A45").Value = 'Internet & Direct Marketing Retail(15)'
$ws.Range("B45").Value = 0.2952069581098644
$ws.Range("A46").Value = 'Health Care Providers & Services(46)'
$ws.Range("B46").Value = 0.2907406161623396
$ws.Range("A47").Value = 'IT Services(52)'
$ws.Range("B47").Value = 0.2895560816946511
$ws.Range("A48").Value = 'Beverages(21)'
$ws.Range("B48").Value = 0.270311358287111
$ws.Range("A50").Value = 'Health Care Equipment & Supplies(83)'
$ws.Range("B50").Value = 0.2530052926919384
$ws.Range("A51").Value = 'Thrifts & Mortgage Finance(47)'
$ws.Range("B51").Value = 0.2454875116583609
$ws.Range("A52").Value = 'Entertainment(22)'
$ws.Range("B52").Value = 0.2426693261706819
$ws.Range("A53").Value = 'Food Products(44)'
$ws.Range("B53").Value = 0.1977162084373913
$ws.Range("A54").Value = 'Pharmaceuticals(48)'
$ws.Range("B54").Value = 0.1616154705497828
$ws.Range("A55").Value = 'Biotechnology(126)'
$ws.Range("B55").Value = 0.1573268451086477

$ws.Range("A56:B59").ClearContents()
